$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A (rows 2-7) values from 1 to 2 ("speed" changed to 2 for both)
$ws.Range("A2:A7").Value = 2

# Update the active selection to A7 to match the saved workbook state
$ws.Range("A7").Select()
